$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for C1, C2, J8, L1, U1 (review corrections / duplicate &
# superseded BOM entries). Deleting from the bottom up keeps earlier row
# numbers valid.
$ws.Rows.Item(24).Delete()   # U1  - IC REG BUCK 3.3V 500MA 10WSON
$ws.Rows.Item(16).Delete()   # L1  - FIXED IND 10UH 1A 540 MOHM SMD
$ws.Rows.Item(15).Delete()   # J8  - TERM BLK 2POS SIDE ENTRY 5MM PCB
$ws.Rows.Item(4).Delete()    # C2  - CAP CER 2.2UF 50V X7R 0805
$ws.Rows.Item(3).Delete()    # C1  - CAP CER 22UF 25V X5R 0805
